$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11.27316440138434
$ws.Range("D2").Value = 16.5658064516129
$ws.Range("E2").Value = 7.177620967741936
$ws.Range("F2").Value = 75.55626456385119
$ws.Range("G2").Value = 6.924969674491517
$ws.Range("H2").Value = 2.792775075844545
$ws.Range("I2").Value = 1023.073500909373
$ws.Range("J2").Value = 44.20833333333334
